# Applies the updated power-flow (vm_pu) results for "case with 380 kV done":
# the slack-bus target voltage (column B) moves from 1.05 to 1.02 p.u. and the
# resulting bus voltages (columns C-F, I-N; column G stays fixed at 1, column H
# has no data) are updated for every data row (rows 2-25) to match the new
# solved values supplied by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ "B" = 1.02; "C" = 1.033209817776665; "D" = 1.040364061206413; "E" = 1.032405580217792; "F" = 1.047653177649391; "I" = 1.02359499962809; "J" = 1.038335797327944; "K" = 1.043146619765433; "L" = 1.035210904742589; "M" = 1.050415214237925; "N" = 1.016343760030791 }
    3 = @{ "B" = 1.02; "C" = 1.035205294500164; "D" = 1.042231984745102; "E" = 1.034140018643152; "F" = 1.049703860060464; "I" = 1.023504579208684; "J" = 1.039968592958495; "K" = 1.04482217148767; "L" = 1.036751634166233; "M" = 1.052274576819547; "N" = 1.01692184552827 }
    4 = @{ "B" = 1.02; "C" = 1.036492677772039; "D" = 1.043437184137004; "E" = 1.035258926093309; "F" = 1.051027453646606; "I" = 1.023443335729026; "J" = 1.041021176213342; "K" = 1.045902485317221; "L" = 1.037744766336921; "M" = 1.053474000180529; "N" = 1.017293601880226 }
    5 = @{ "B" = 1.02; "C" = 1.037033002415418; "D" = 1.043943040265627; "E" = 1.035728522832348; "F" = 1.051583117563117; "I" = 1.023416932628352; "J" = 1.041462757539928; "K" = 1.046355740361174; "L" = 1.038161382251197; "M" = 1.053977371825857; "N" = 1.017449343421187 }
    6 = @{ "B" = 1.02; "C" = 1.037123673732544; "D" = 1.044027928877299; "E" = 1.035807324319759; "F" = 1.051676371226323; "I" = 1.023412460900047; "J" = 1.041536847334351; "K" = 1.046431791159426; "L" = 1.038231281835209; "M" = 1.054061839995004; "N" = 1.017475461354829 }
    7 = @{ "B" = 1.02; "C" = 1.036499901085955; "D" = 1.043443946563794; "E" = 1.035265203953129; "F" = 1.051034881467698; "I" = 1.02344298551034; "J" = 1.041027080251329; "K" = 1.045908545277516; "L" = 1.037750336678346; "M" = 1.053480729630006; "N" = 1.017295685038061 }
    8 = @{ "B" = 1.02; "C" = 1.033885004760162; "D" = 1.040996066642931; "E" = 1.032992455666477; "F" = 1.048346920193924; "I" = 1.023565007130014; "J" = 1.038888436951764; "K" = 1.043713695853899; "L" = 1.035732404095281; "M" = 1.051044375663712; "N" = 1.016539607745481 }
    9 = @{ "B" = 1.02; "C" = 1.029246850724944; "D" = 1.036655025604318; "E" = 1.028960729416596; "F" = 1.0435837807031; "I" = 1.023759139049229; "J" = 1.035088772547041; "K" = 1.039815468083833; "L" = 1.032146442497025; "M" = 1.046721816072014; "N" = 1.015189358393274 }
    10 = @{ "B" = 1.02; "C" = 1.026132796579068; "D" = 1.033741098287959; "E" = 1.026253606906328; "F" = 1.040388965633012; "I" = 1.023874603872166; "J" = 1.0325334997304; "K" = 1.037194787113908; "L" = 1.029734393480451; "M" = 1.043818955894509; "N" = 1.0142767017768 }
    11 = @{ "B" = 1.02; "C" = 1.024778827063529; "D" = 1.032474308917279; "E" = 1.025076533615259; "F" = 1.039000640693447; "I" = 1.023921309042974; "J" = 1.031421500692188; "K" = 1.036054532196016; "L" = 1.028684610337648; "M" = 1.042556655902951; "N" = 1.013878453048855 }
    12 = @{ "B" = 1.02; "C" = 1.024275036930368; "D" = 1.032002983224955; "E" = 1.024638560003959; "F" = 1.038484182377525; "I" = 1.023938164268129; "J" = 1.031007596704836; "K" = 1.035630142264607; "L" = 1.028293847671017; "M" = 1.042086951809284; "N" = 1.013730057304736 }
    13 = @{ "B" = 1.02; "C" = 1.024383141274503; "D" = 1.032104120081374; "E" = 1.024732541443186; "F" = 1.038594999873627; "I" = 1.023934571058996; "J" = 1.031096419737626; "K" = 1.035721214142586; "L" = 1.028377705356544; "M" = 1.042187742914528; "N" = 1.013761910028578 }
    14 = @{ "B" = 1.02; "C" = 1.024737201444767; "D" = 1.032435365131201; "E" = 1.025040346147723; "F" = 1.038957966018964; "I" = 1.023922712353276; "J" = 1.031387304905238; "K" = 1.036019469445282; "L" = 1.028652326765827; "M" = 1.042517847098966; "N" = 1.01386619621502 }
    15 = @{ "B" = 1.02; "C" = 1.024955234008424; "D" = 1.032639351487304; "E" = 1.025229893850178; "F" = 1.039181498183369; "I" = 1.023915340507297; "J" = 1.031566414425138; "K" = 1.036203121201448; "L" = 1.028821419880529; "M" = 1.042721124426384; "N" = 1.013930388052359 }
    16 = @{ "B" = 1.02; "C" = 1.026222535333187; "D" = 1.033825062488291; "E" = 1.026331620682546; "F" = 1.040480997444552; "I" = 1.023871434968932; "J" = 1.032607180461786; "K" = 1.037270344456556; "L" = 1.029803949482805; "M" = 1.043902615730691; "N" = 1.014303067002877 }
    17 = @{ "B" = 1.02; "C" = 1.027015969617862; "D" = 1.034567459787456; "E" = 1.027021383557602; "F" = 1.041294793619089; "I" = 1.023843013862359; "J" = 1.033258522782283; "K" = 1.037938300148852; "L" = 1.030418816083113; "M" = 1.044642284737378; "N" = 1.014536012922803 }
    18 = @{ "B" = 1.02; "C" = 1.027478231641089; "D" = 1.03500000313566; "E" = 1.027423242214803; "F" = 1.041768991645075; "I" = 1.023826118601784; "J" = 1.033637905997685; "K" = 1.03832738008622; "L" = 1.030776942719599; "M" = 1.045073207538318; "N" = 1.014671591430781 }
    19 = @{ "B" = 1.02; "C" = 1.027635761231493; "D" = 1.035147407869116; "E" = 1.02756018692078; "F" = 1.041930601136433; "I" = 1.023820303822167; "J" = 1.033767176070035; "K" = 1.038459957592012; "L" = 1.030898968036121; "M" = 1.045220054771652; "N" = 1.014717770427696 }
    20 = @{ "B" = 1.02; "C" = 1.026930897155922; "D" = 1.034487857882049; "E" = 1.026947427148485; "F" = 1.04120753038543; "I" = 1.023846096030143; "J" = 1.033188695298134; "K" = 1.037866689545705; "L" = 1.030352900111843; "M" = 1.044562978573577; "N" = 1.014511050593771 }
    21 = @{ "B" = 1.02; "C" = 1.02463296368405; "D" = 1.032337843569791; "E" = 1.024949726405382; "F" = 1.038851103099653; "I" = 1.02392621804881; "J" = 1.031301670356717; "K" = 1.035931664315433; "L" = 1.028571480603898; "M" = 1.04242066276641; "N" = 1.013835499531121 }
    22 = @{ "B" = 1.02; "C" = 1.023183138755013; "D" = 1.030981496454732; "E" = 1.023689305324463; "F" = 1.037365038888681; "I" = 1.023973741382756; "J" = 1.030110246246444; "K" = 1.034710115551421; "L" = 1.027446638298702; "M" = 1.041068889878637; "N" = 1.01340803931967 }
    23 = @{ "B" = 1.02; "C" = 1.023952204531718; "D" = 1.031700961772326; "E" = 1.024357902465237; "F" = 1.038153264580367; "I" = 1.023948818288664; "J" = 1.030742322659691; "K" = 1.035358156535574; "L" = 1.028043400484504; "M" = 1.04178595569624; "N" = 1.013634904183486 }
    24 = @{ "B" = 1.02; "C" = 1.026969339391033; "D" = 1.034523828052759; "E" = 1.026980846312841; "F" = 1.041246962352938; "I" = 1.023844704313072; "J" = 1.033220248981015; "K" = 1.03789904892345; "L" = 1.037744766336921; "M" = 1.044598815201521; "N" = 1.014522330906809 }
    25 = @{ "B" = 1.02; "C" = 1.030449680348662; "D" = 1.037780694888123; "E" = 1.030006337236424; "F" = 1.044818472727445; "I" = 1.02371141742819; "J" = 1.036074890452035; "K" = 1.040827012393664; "L" = 1.033077186371113; "M" = 1.015540600153758 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value = $data[$row][$col]
    }
}
